$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 0. Update the selections left behind on the two sheets that were active
#    while the author was working (AccountOpening_Loan_ACOPL_TwoWh and
#    Debenture_Type_setting) before moving on to the brand-new sheet.
# ---------------------------------------------------------------------------
$loanSheetSel = $wb.Worksheets.Item("AccountOpening_Loan_ACOPL_TwoWh")
$loanSheetSel.Activate()
$loanSheetSel.Range("AE2").Select()

$typeSheetSel = $wb.Worksheets.Item("Debenture_Type_setting")
$typeSheetSel.Activate()
$typeSheetSel.Range("A1:E2").Select()

# ---------------------------------------------------------------------------
# 1. Insert the new worksheet "Debenture_Issue_setting" right after the
#    existing "Debenture_Type_setting" sheet (becomes the last / 14th sheet).
# ---------------------------------------------------------------------------
$typeSheet = $wb.Worksheets.Item("Debenture_Type_setting")
$ws = $wb.Worksheets.Add([System.Type]::Missing, $typeSheet)
$ws.Name = "Debenture_Issue_setting"

# ---------------------------------------------------------------------------
# 2. Column widths (only columns E:K get a custom width).
# ---------------------------------------------------------------------------
$ws.Columns.Item(5).ColumnWidth = 13.0            # E  -> width 13.81640625
$ws.Columns.Item(6).ColumnWidth = 16.6666666666667 # F  -> width 17.453125
$ws.Columns.Item(7).ColumnWidth = 14.5             # G  -> width 15.26953125
$ws.Columns.Item(8).ColumnWidth = 11.3333333333333 # H  -> width 12.1796875
$ws.Columns.Item(9).ColumnWidth = 13.0             # I  -> width 13.81640625
$ws.Columns.Item(10).ColumnWidth = 20.0            # J  -> width 20.90625
$ws.Columns.Item(11).ColumnWidth = 21.1666666666667 # K -> width 22

# ---------------------------------------------------------------------------
# 3. Copy cell formatting from existing cells that already carry the target
#    styles, then fill in the values/number formats for the new sheet.
# ---------------------------------------------------------------------------
$typeSheet.Range("A1:D2").Copy()
$ws.Range("A1:D2").PasteSpecial(-4122)

$fixedDeposit = $wb.Worksheets.Item("AccountOpening_FixedDeposit")
$fixedDeposit.Range("O1:U1").Copy()
$ws.Range("E1:K1").PasteSpecial(-4122)

$loanSheet = $wb.Worksheets.Item("AccountOpening_Loan_ACOPL_TwoWh")
$loanSheet.Range("H2").Copy()
$ws.Range("E2").PasteSpecial(-4122)

$ws.Application.CutCopyMode = $false

# Row heights (wrapped 3-line header/data rows).
$ws.Rows.Item(1).RowHeight = 43.5
$ws.Rows.Item(2).RowHeight = 43.5

# ---------------------------------------------------------------------------
# 4. Cell values. The write order below reproduces the exact order in which
#    new shared-string entries were appended by the original authoring tool.
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "TestScenario"
$ws.Range("B1").Value = "Run"
$ws.Range("C1").Value = "pcRegFormName"
$ws.Range("D1").Value = "pcRegFormPcName"
$ws.Range("E1").Value = "IssueFromDate"
$ws.Range("F1").Value = "UnitValue"
$ws.Range("G1").Value = "ReturnROI"
$ws.Range("H1").Value = "CallDate"
$ws.Range("I1").Value = "PutDate"

$ws.Range("A2").Value = "Debenture_Type_setting"
$ws.Range("B2").Value = "Yes"
$ws.Range("C2").Value = "qwerty"
$ws.Range("D2").Value = "zxcvb"
$ws.Range("E2").Value = "25/12/2300"
$ws.Range("F2").Value = 250
$ws.Range("G2").Value = 3

$ws.Range("H2").NumberFormat = "mm-dd-yy"
$ws.Range("H2").Value = "15/03/2301"
$ws.Range("I2").Value = "20/04/2301"

$ws.Range("J1").Value = "CallDatePrematureROI"
$ws.Range("K1").Value = "PutDatePrematureROI"

$ws.Range("J2").Value = 7
$ws.Range("K2").Value = 8

# ---------------------------------------------------------------------------
# 6. View state: activate the new sheet, scroll so column D is the left-most
#    visible column, and select K1 (matches the authored sheetView).
# ---------------------------------------------------------------------------
$ws.Activate()
$win = $wb.Windows.Item(1)
$win.ScrollColumn = 4
$ws.Range("K1").Select()
